$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text-formatted numeric-looking price strings (e.g. "28.534.52").
# Force the Price column to stay Text so Excel does not renormalize values
# (stripping trailing zeros, switching to scientific notation, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$changes = @(
    @(2, 4, '28.534.52'),
    @(2, 5, '  +0.49%  '),
    @(3, 4, '1.870.82'),
    @(3, 5, '  -0.47%  '),
    @(4, 5, '  -0.23%  '),
    @(5, 4, '314.93'),
    @(5, 5, '  -0.23%  '),
    @(6, 4, '1.007'),
    @(6, 5, '  -0.51%  '),
    @(7, 4, '0.5062'),
    @(7, 5, '  -1.50%  '),
    @(8, 4, '0.3893'),
    @(8, 5, '  -1.09%  '),
    @(9, 4, '0.08351'),
    @(9, 5, '  +0.37%  '),
    @(10, 2, 'OKB'),
    @(10, 3, 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'),
    @(10, 4, '41.75'),
    @(10, 5, '  -0.42%  '),
    @(11, 2, 'Polygon'),
    @(11, 3, 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'),
    @(11, 4, '1.101'),
    @(11, 5, '  -1.88%  '),
    @(12, 2, 'Polkadot'),
    @(12, 3, 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'),
    @(12, 4, '6.204'),
    @(12, 5, '  -1.24%  '),
    @(13, 2, 'WrappedEther'),
    @(13, 3, 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'),
    @(13, 4, '1.871.22'),
    @(13, 5, '  -0.13%  '),
    @(14, 2, 'Solana'),
    @(14, 3, 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'),
    @(14, 4, '20.36'),
    @(14, 5, '  -0.22%  '),
    @(15, 2, 'Chainlink'),
    @(15, 3, 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'),
    @(15, 4, '7.233'),
    @(15, 5, '  -0.30%  '),
    @(16, 2, 'BinanceUSD'),
    @(16, 3, 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'),
    @(16, 4, '1.008'),
    @(16, 5, '  -0.29%  '),
    @(17, 2, 'ShibaInu'),
    @(17, 3, 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'),
    @(17, 4, '0.00001101'),
    @(17, 5, '  -0.61%  '),
    @(18, 2, 'Litecoin'),
    @(18, 3, 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'),
    @(18, 4, '90.94'),
    @(18, 5, '  -0.47%  '),
    @(19, 2, 'TRON'),
    @(19, 3, 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'),
    @(19, 4, '0.06697'),
    @(19, 5, '  -0.63%  '),
    @(20, 2, 'Avalanche'),
    @(20, 3, 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'),
    @(20, 4, '17.64'),
    @(20, 5, '  -0.65%  '),
    @(21, 2, 'Dai'),
    @(21, 3, 'https://coinranking.com/coin/MoTuySvg7+dai-dai'),
    @(21, 4, '1.007'),
    @(21, 5, '  -0.58%  '),
    @(22, 2, 'Uniswap'),
    @(22, 3, 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'),
    @(22, 4, '5.913'),
    @(22, 5, '  -1.62%  '),
    @(23, 2, 'WrappedBTC'),
    @(23, 3, 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'),
    @(23, 4, '28.571.43'),
    @(23, 5, '  +0.47%  '),
    @(24, 2, 'Cosmos'),
    @(24, 3, 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'),
    @(24, 4, '11.04'),
    @(24, 5, '  -1.12%  '),
    @(25, 2, 'Toncoin'),
    @(25, 3, 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'),
    @(25, 4, '2.233'),
    @(25, 5, '  -0.64%  '),
    @(26, 2, 'WrappedliquidstakedEther2.0'),
    @(26, 3, 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'),
    @(26, 4, '2.086.29'),
    @(26, 5, '  -0.16%  '),
    @(27, 2, 'Monero'),
    @(27, 3, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'),
    @(27, 4, '161.61'),
    @(27, 5, '  +0.51%  '),
    @(28, 2, 'EthereumClassic'),
    @(28, 3, 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'),
    @(28, 4, '20.62'),
    @(28, 5, '  -0.98%  '),
    @(29, 2, 'LidoDAOToken'),
    @(29, 3, 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'),
    @(29, 4, '2.338'),
    @(29, 5, '  -4.52%  '),
    @(30, 2, 'BitcoinCash'),
    @(30, 3, 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'),
    @(30, 4, '125.76'),
    @(30, 5, '  -1.11%  '),
    @(31, 2, 'Stellar'),
    @(31, 3, 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'),
    @(31, 4, '0.1041'),
    @(31, 5, '  -2.31%  '),
    @(32, 2, 'ImmutableX'),
    @(32, 3, 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @(32, 4, '1.039'),
    @(32, 5, '  -1.22%  '),
    @(33, 2, 'Filecoin'),
    @(33, 3, 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @(33, 4, '5.777'),
    @(33, 5, '  -2.11%  '),
    @(34, 2, 'HuobiToken'),
    @(34, 3, 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'),
    @(34, 4, '3.612'),
    @(34, 5, '  -0.74%  '),
    @(35, 2, 'VeChain'),
    @(35, 3, 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'),
    @(35, 4, '0.02446'),
    @(35, 5, '  +0.04%  '),
    @(36, 2, 'Hedera'),
    @(36, 3, 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @(36, 4, '0.06528'),
    @(36, 5, '  -0.09%  '),
    @(37, 2, 'Algorand'),
    @(37, 3, 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'),
    @(37, 4, '0.2158'),
    @(37, 5, '  -1.53%  '),
    @(38, 2, 'FraxShare'),
    @(38, 3, 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'),
    @(38, 4, '8.840'),
    @(38, 5, '  -4.19%  '),
    @(39, 2, 'InternetComputer(DFINITY)'),
    @(39, 3, 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'),
    @(39, 4, '5.045'),
    @(39, 5, '  +1.04%  '),
    @(40, 2, 'TrustWalletToken'),
    @(40, 3, 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'),
    @(40, 4, '1.254'),
    @(40, 5, '  -0.34%  '),
    @(41, 2, 'ARBITRUM'),
    @(41, 3, 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'),
    @(41, 4, '1.186'),
    @(41, 5, '  -0.46%  '),
    @(42, 2, 'TheSandbox'),
    @(42, 3, 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'),
    @(42, 4, '0.6417'),
    @(42, 5, '  -1.05%  '),
    @(43, 2, 'Aptos'),
    @(43, 3, 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'),
    @(43, 4, '11.08'),
    @(43, 5, '  -1.08%  '),
    @(44, 2, 'Frax'),
    @(44, 3, 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'),
    @(44, 4, '1.007'),
    @(44, 5, '  -0.69%  '),
    @(45, 2, 'Decentraland'),
    @(45, 3, 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'),
    @(45, 4, '0.6012'),
    @(45, 5, '  -0.98%  '),
    @(46, 2, 'EnergySwap'),
    @(46, 3, 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'),
    @(46, 4, '12.95'),
    @(46, 5, '  -1.70%  '),
    @(47, 2, 'PancakeSwap'),
    @(47, 3, 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'),
    @(47, 4, '3.688'),
    @(47, 5, '  -0.45%  '),
    @(48, 2, 'NEARProtocol'),
    @(48, 3, 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'),
    @(48, 4, '2.000'),
    @(48, 5, '  -1.27%  '),
    @(49, 2, 'EOS'),
    @(49, 3, 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'),
    @(49, 4, '1.212'),
    @(49, 5, '  -0.86%  '),
    @(50, 2, 'Quant'),
    @(50, 3, 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'),
    @(50, 4, '121.79'),
    @(50, 5, '  -0.22%  '),
    @(51, 2, 'WEMIXTOKEN'),
    @(51, 3, 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'),
    @(51, 4, '1.178'),
    @(51, 5, '  -8.09%  ')
)

foreach ($item in $changes) {
    $r = $item[0]
    $c = $item[1]
    $v = $item[2]
    $ws.Cells.Item($r, $c).Value = $v
}

Write-Output "Applied $($changes.Count) cell updates"
